$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A11").Value = "bunt.saar"
$ws.Range("B11").Value = "bunt.saar"
$ws.Range("H11").Value = "https://de.wikipedia.org/wiki/Bunt.saar"
$ws.Range("B12").Value = "SGV Solidarität, Gerechtigkeit, Veränderung"
$ws.Range("A12").Value = "SGV"
$ws.Range("B12").Select() | Out-Null
